$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch P51 (a column-P cell) so the row block 49:64 recomputes its
# "spans" metadata out to column P, then clear it back to blank -
# matches the target's <c r="P51"/> placeholder cell.
$ws.Range("P51").Value = 1
$ws.Range("P51").ClearContents()

# Attendance/score updates for rows 53-70 (columns M/N)
$ws.Range("M53").Value = 11

$ws.Range("M54").Value = 8

$ws.Range("M55").Value = 7

$ws.Range("M56").Value = 15

$ws.Range("M57").ClearContents()

$ws.Range("M58").Value = 9
$ws.Range("N58").Value = 8

$ws.Range("N59").Value = 7

$ws.Range("M60").Value = 6
$ws.Range("N60").Value = 1

$ws.Range("M61").Value = 10
$ws.Range("N61").Value = 17

$ws.Range("M62").Value = 17
$ws.Range("N62").Value = 16

$ws.Range("M63").Value = 14
$ws.Range("N63").Value = 13

$ws.Range("M64").Value = 13
$ws.Range("N64").Value = 4

$ws.Range("M65").Value = 5
$ws.Range("N65").Value = 3

$ws.Range("M66").Value = 2
$ws.Range("N66").Value = 2

$ws.Range("M67").Value = 3
$ws.Range("N67").Value = 9

$ws.Range("M68").Value = 16
$ws.Range("N68").Value = 15

$ws.Range("M69").Value = 15
$ws.Range("N69").Value = 14

$ws.Range("M70").Value = 18
$ws.Range("N70").Value = 15

# Restore the view: scroll so row 43 is at the top and select N58
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("N58").Select()
